$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column C (date "Förändrad" column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 493 }

$rng = $ws.Range("C2:C$lastRow")
foreach ($cell in $rng.Cells) {
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
